# Increase the used range by one column (A1:G6 -> A1:H6) and populate the
# new "maatschappij" column H with its header ("P368") and the value
# for the second data row ("x"). Excel will automatically grow the
# worksheet's <dimension> to cover the newly written cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "P368"
$ws.Range("H2").Value = "x"
